$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # RUNMANAGER
$ws2 = $wb.Worksheets.Item(2)   # DATA

# --- DATA sheet: update browser "version" column values -----------------
# Existing values: "88.0.4324.96" -> "94.0.4606.61"
#                  "69.0"         -> "92.0.1"
# (the password/base64 column values are unaffected; only their shared
# string slot shuffles as a side effect of the sst edit)
# Set the "92.0.1" replacements first so the new shared-string entries land
# in the same order as the target workbook.
$ws2.Range("D3").Value = "'92.0.1"
$ws2.Range("D5").Value = "'92.0.1"
$ws2.Range("D9").Value = "'92.0.1"

$ws2.Range("D2").Value = "'94.0.4606.61"
$ws2.Range("D4").Value = "'94.0.4606.61"
$ws2.Range("D8").Value = "'94.0.4606.61"

# --- DATA sheet: column D width tweak (auto bestFit -> fixed width) -----
$ws2.Columns.Item(4).ColumnWidth = 24.83

# --- Selection / active cell updates -------------------------------------
[void]$ws1.Range("D2").Select()
[void]$ws2.Range("D8").Select()
